$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)
$ws.Activate()

# Widen column F slightly (target stored width ~5.75 chars; engine quantizes to 1/7 steps)
$ws.Columns.Item(6).ColumnWidth = 4.93

# Row 23
$ws.Range("E23").Value = 4
$ws.Range("F23").Value = "done"
$ws.Range("J23").Value = "21:00-23:00"
$ws.Range("L23").Value = "2h"

# Row 24
$ws.Range("J24").Value = "23:00-23:30"
$ws.Range("L24").Value = "2.5h"

# Row 25
$ws.Range("F25").Value = "a"
$ws.Range("G25").Value = "jquery"
$ws.Range("L25").Value = "study and integrate"

# Row 26
$ws.Range("F26").Value = "b"
$ws.Range("G26").Value = "ui"

# Row 27-35 (jquery ui widgets studied)
$ws.Range("H27").Value = "Button"
$ws.Range("H28").Value = "Datepicker"
$ws.Range("H29").Value = "Autocomplete"
$ws.Range("H30").Value = "Progressbar"
$ws.Range("H31").Value = "Dialog"
$ws.Range("I31").Value = "model dialog, model form"
$ws.Range("H32").Value = "Tabs"
$ws.Range("H33").Value = "ToggleClass"
$ws.Range("H34").Value = "Effect"
$ws.Range("H35").Value = "Position"

# Row 37
$ws.Range("F37").Value = "c"
$ws.Range("G37").Value = "theme"

# Row 38
$ws.Range("F38").Value = "d"
$ws.Range("G38").Value = "css framework"

# Row 39
$ws.Range("F39").Value = "e"
$ws.Range("G39").Value = "layout framework"

# Update selection to match target view state
$ws.Range("L29").Select()
